$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'City Planning Commission Approves Jamaica Neighborhood Plan In Queens'
$ws.Cells.Item(2, 2).Value = 'https://newyorkyimby.com/2025/08/city-planning-commission-approves-jamaica-neighborhood-plan-in-queens.html'
$ws.Cells.Item(2, 3).Value = 'The New York City Planning Commission has approved the Jamaica Neighborhood Plan, a rezoning initiative aimed at transforming 230 blocks of <a href="https://newyorkyimby.com/neighborhoods/jamaica">Jamaica</a>, Queens. Led by Mayor Eric Adams and the NYC Department of City Planning, the plan is expected to create over 12,000 new homes, including 4,000 permanently affordable units, along with 7,000 new jobs and more than 2 million square feet of commercial and community space. The plan will now advance to the City Council for final review after two years of development with extensive community engagement.'
$ws.Cells.Item(2, 4).Value = '2025-08-14T11:30:52+00:00'
$ws.Cells.Item(2, 5).Value = 'Thu, 14 Aug 2025 11:30:52 +0000'
$ws.Cells.Item(2, 6).Value = 'YIMBY'
$ws.Cells.Item(2, 7).Value = 'YIMBY - Jamaica'

# Row 3
$ws.Cells.Item(3, 1).Value = 'Permits Filed for 64-23 Austin Street in Rego Park, Queens'
$ws.Cells.Item(3, 2).Value = 'https://newyorkyimby.com/2025/06/permits-filed-for-64-23-austin-street-in-rego-park-queens.html'
$ws.Cells.Item(3, 3).Value = 'Permits have been filed for a six-story residential building at 64-23 Austin Street in <a href="https://newyorkyimby.com/neighborhoods/rego-park">Rego Park</a>, Queens. Located between 63rd Drive and 64th Road, the lot is near the 63 Drive-Rego Park subway station, served by the M and R trains. Danil Ilyayev is listed as the owner behind the applications.'
$ws.Cells.Item(3, 4).Value = '2025-06-19T10:30:09+00:00'
$ws.Cells.Item(3, 5).Value = 'Thu, 19 Jun 2025 10:30:09 +0000'
$ws.Cells.Item(3, 6).Value = 'YIMBY'
$ws.Cells.Item(3, 7).Value = 'YIMBY - Rego Park'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Demolition Nears Completion at 54-03, 54-11 Queens Boulevard in Woodside, Queens'
$ws.Cells.Item(4, 2).Value = 'https://newyorkyimby.com/2025/03/demolition-nears-completion-at-54-03-54-11-queens-boulevard-in-woodside-queens.html'
$ws.Cells.Item(4, 3).Value = 'Demolition is almost done at 54-03 and 54-11 Queens Boulevard in <a href="https://newyorkyimby.com/neighborhoods/woodside">Woodside</a>, Queens. The rectangular 12,000-square-foot property is bound by Queens Boulevard to the south, 55th Street to the east, and 54th Street to the west. Santiago Helman of Nostrand Land LLC is listed as the owner and filed demolition permits, in part, due to the instability of both structures.'
$ws.Cells.Item(4, 4).Value = '2025-03-23T11:30:59+00:00'
$ws.Cells.Item(4, 5).Value = 'Sun, 23 Mar 2025 11:30:59 +0000'
$ws.Cells.Item(4, 6).Value = 'YIMBY'
$ws.Cells.Item(4, 7).Value = 'YIMBY - Woodside'

# Row 5
$ws.Cells.Item(5, 1).Value = '38-15 Queens Boulevard’s Steel Superstructure Rises in Sunnyside, Queens'
$ws.Cells.Item(5, 2).Value = 'https://newyorkyimby.com/2024/11/38-15-queens-boulevards-steel-superstructure-rises-in-sunnyside-queens.html'
$ws.Cells.Item(5, 3).Value = 'Construction is rising on 38-15 Queens Boulevard, a two-story car dealership in <a href="https://newyorkyimby.com/neighborhoods/sunnyside">Sunnyside</a>, Queens. Developed by Auto Group, which purchased the property along with the abutting 38-08 43rd Avenue for $38 million under the DCO QOZB 38-15 QB LLC in early 2023, the structure will cover 0.69 acres and contain an expansive auto showroom and an ajoining service facility. The project is bound by Queens Boulevard to the south and 39th Street to the east.'
$ws.Cells.Item(5, 4).Value = '2024-11-10T12:30:35+00:00'
$ws.Cells.Item(5, 5).Value = 'Sun, 10 Nov 2024 12:30:35 +0000'
$ws.Cells.Item(5, 6).Value = 'YIMBY'
$ws.Cells.Item(5, 7).Value = 'YIMBY - Sunnyside'

# Row 6
$ws.Cells.Item(6, 1).Value = 'Developers Secure Refinancing for The Yellowstone at 69-65 Yellowstone Boulevard in Forest Hills, Queens'
$ws.Cells.Item(6, 2).Value = 'https://newyorkyimby.com/2024/03/developers-secure-refinancing-for-the-yellowstone-at-69-65-yellowstone-boulevard-in-forest-hills-queens.html'
$ws.Cells.Item(6, 3).Value = '<a href="https://newyorkyimby.com/category/slate-property-group">Slate Property Group</a>, in collaboration with Grobman Gross Properties, has finalized a $97 million refinancing deal for The Yellowstone, a recently completed 11-story residential building at <a href="https://newyorkyimby.com/category/69-65-yellowstone-boulevard">69-65 Yellowstone Boulevard</a> in <a href="https://newyorkyimby.com/neighborhoods/forest-hills">Forest Hills</a>, <a href="https://newyorkyimby.com/category/queens">Queens</a>. The refinancing, which will be provided by Apollo Global Real Estate Management L.P., will be used to settle the building''s existing debts and support its financial structure.'
$ws.Cells.Item(6, 4).Value = '2024-03-05T12:00:48+00:00'
$ws.Cells.Item(6, 5).Value = 'Tue, 05 Mar 2024 12:00:48 +0000'
$ws.Cells.Item(6, 6).Value = 'YIMBY'
$ws.Cells.Item(6, 7).Value = 'YIMBY - Forest Hills'

# Row 7
$ws.Cells.Item(7, 1).Value = 'RYBAK Development Secures $71M Loan for The Austin Condominiums in Forest Hills, Queens'
$ws.Cells.Item(7, 2).Value = 'https://newyorkyimby.com/2024/01/rybak-development-secures-71m-loan-for-the-austin-condominiums-in-forest-hills-queens.html'
$ws.Cells.Item(7, 3).Value = '<a href="https://newyorkyimby.com/category/rybak-development">RYBAK Development</a> has successfully secured a $71 million construction loan for The Austin, a new condominium project located at <a href="https://newyorkyimby.com/category/78-29-austin-street">78-29 Austin Street</a> in <a href="https://newyorkyimby.com/neighborhoods/forest-hills">Forest Hills</a>, <a href="https://newyorkyimby.com/category/queens">Queens</a>. <a href="https://newyorkyimby.com/category/jll-capital-markets">JLL Capital Markets</a> facilitated the financing for the project from Valley National Bank.'
$ws.Cells.Item(7, 4).Value = '2024-01-15T12:00:09+00:00'
$ws.Cells.Item(7, 5).Value = 'Mon, 15 Jan 2024 12:00:09 +0000'
$ws.Cells.Item(7, 6).Value = 'YIMBY'
$ws.Cells.Item(7, 7).Value = 'YIMBY - Forest Hills'

# Row 8
$ws.Cells.Item(8, 1).Value = 'Permits Filed for 48-37 48th Street in Sunnyside, Queens'
$ws.Cells.Item(8, 2).Value = 'https://newyorkyimby.com/2023/07/permits-filed-for-48-37-48th-street-in-sunnyside-queens.html'
$ws.Cells.Item(8, 3).Value = 'Permits have been filed for a five-story commercial building at 48-37 48th Street in <a href="https://newyorkyimby.com/category/sunnyside">Sunnyside</a>, Queens. Located between 48th Avenue and 50th Avenue, the lot is near the 52nd Street subway station, serviced by the 7 train. Denis Iserovich of <a href="https://newyorkyimby.com/category/east-end-capital">East End Capital</a> is listed as the owner behind the applications.'
$ws.Cells.Item(8, 4).Value = '2023-07-14T10:30:24+00:00'
$ws.Cells.Item(8, 5).Value = 'Fri, 14 Jul 2023 10:30:24 +0000'
$ws.Cells.Item(8, 6).Value = 'YIMBY'
$ws.Cells.Item(8, 7).Value = 'YIMBY - Sunnyside'

# Row 9
$ws.Cells.Item(9, 1).Value = 'Myrtle Point Rental Tower Tops Out at 3-50 St. Nicholas Avenue in Ridgewood, Queens'
$ws.Cells.Item(9, 2).Value = 'https://newyorkyimby.com/2022/12/myrtle-point-rental-tower-tops-out-at-3-50-st-nicholas-avenue-in-ridgewood-queens.html'
$ws.Cells.Item(9, 3).Value = 'Construction has topped out on Myrtle Point, a 17-story rental tower at <a href="https://newyorkyimby.com/category/3-50-st-nicholas-avenue-2">3-50 St. Nicholas Avenue</a> in <a href="https://newyorkyimby.com/neighborhoods/ridgewood">Ridgewood</a>, Queens. Designed by<a href="https://newyorkyimby.com/category/s9-architecture-engineering"> S9 Architecture</a> and developed in collaboration by Arch Companies and <a href="https://newyorkyimby.com/category/ab-capstone">AB Capstone</a>, the building will be the tallest residential property in the Ridgewood-<a href="https://newyorkyimby.com/neighborhoods/bushwick">Bushwick</a> area and is expected to debut in 2023.'
$ws.Cells.Item(9, 4).Value = '2022-12-05T12:00:04+00:00'
$ws.Cells.Item(9, 5).Value = 'Mon, 05 Dec 2022 12:00:04 +0000'
$ws.Cells.Item(9, 6).Value = 'YIMBY'
$ws.Cells.Item(9, 7).Value = 'YIMBY - Ridgewood'

